$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.24"
$ws.Range("E2").Value = "'2.28%"
$ws.Range("D3").Value = "'43.21"
$ws.Range("E3").Value = "'6.48%"
$ws.Range("D4").Value = "'5.088"
$ws.Range("E4").Value = "'1.42%"
$ws.Range("D5").Value = "'0.07663"
$ws.Range("E5").Value = "'3.16%"
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = "'4.410"
$ws.Range("E6").Value = "'1.99%"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = "'1.614"
$ws.Range("E7").Value = "'2.66%"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'1.006"
$ws.Range("E8").Value = "'8.04%"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = "'0.1243"
$ws.Range("E9").Value = "'4.17%"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1858"
$ws.Range("E10").Value = "'3.41%"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.09115"
$ws.Range("E11").Value = "'4.01%"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.04168"
$ws.Range("E12").Value = "'-2.47%"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.1047"
$ws.Range("E13").Value = "'-0.55%"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001295"
$ws.Range("E14").Value = "'2.35%"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = "'0.005743"
$ws.Range("E15").Value = "'-3.44%"
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = "'0.007430"
$ws.Range("E16").Value = "'1,900.86%"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.334"
$ws.Range("E17").Value = "'-0.24%"
$ws.Range("D18").Value = "'2.355"
$ws.Range("E18").Value = "'-1.62%"
$ws.Range("D19").Value = "'0.3353"
$ws.Range("E19").Value = "'1.70%"
$ws.Range("D20").Value = "'8.423"
$ws.Range("E20").Value = "'7.07%"
$ws.Range("D21").Value = "'0.1397"
$ws.Range("E21").Value = "'1.33%"
$ws.Range("D22").Value = "'0.3193"
$ws.Range("D23").Value = "'0.04174"
$ws.Range("E23").Value = "'5.75%"
$ws.Range("D24").Value = "'0.001286"
$ws.Range("E24").Value = "'1.59%"
$ws.Range("D25").Value = "'0.004498"
$ws.Range("E25").Value = "'17.69%"
$ws.Range("D26").Value = "'0.0001352"
$ws.Range("E26").Value = "'10.66%"
$ws.Range("D38").Value = "'0.02460"
$ws.Range("E38").Value = "'4.22%"
$ws.Range("D39").Value = "'0.05286"
$ws.Range("E39").Value = "'3.18%"
$ws.Range("D40").Value = "'0.005975"
$ws.Range("E40").Value = "'1.27%"
$ws.Range("D41").Value = "'0.007676"
$ws.Range("E41").Value = "'-0.66%"
$ws.Range("D42").Value = "'0.1348"
$ws.Range("E42").Value = "'2.35%"
$ws.Range("D43").Value = "'0.007366"
$ws.Range("E43").Value = "'-0.06%"
$ws.Range("D44").Value = "'0.007553"
$ws.Range("E44").Value = "'7.32%"
$ws.Range("D45").Value = "'0.3029"
$ws.Range("E45").Value = "'3.47%"
$ws.Range("D46").Value = "'0.00006697"
$ws.Range("E46").Value = "'7.92%"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("D48").Value = "'0.04497"
$ws.Range("E48").Value = "'-2.88%"
$ws.Range("E49").Value = "'0.24%"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("E51").Value = "'0.04%"
